# Update BAU CCS amounts for latest BAU trajectory

$wb = $excel.ActiveWorkbook

# 1. Replace the ": NoSettings" suffix with ": test" across the row labels
#    in the "BAU Emissions" sheet (column A text entries referencing the
#    shared-string table).
$wsBAU = $wb.Worksheets.Item("BAU Emissions")
$replaceResult = $wsBAU.UsedRange.Replace(" : NoSettings", " : test")

# 2. Update the BAU trajectory amounts for the "Industrial Sector Energy
#    Related Emissions before CCS[natural gas if,iron and steel 241,CO2]"
#    row (row 94), columns M (2032) through AE (2050).
$newValues = @(1001080, 2002150, 3003230, 4004300, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380)
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $wsBAU.Cells.Item(94, 13 + $i).Value2 = $newValues[$i]
}

# 3. Update the view state: select A30:AE280 on the "BAU Emissions" sheet,
#    then make "About" the active sheet/tab.
$activateResult = $wsBAU.Activate()
$selectResult = $wsBAU.Range("A30:AE280").Select()

$wsAbout = $wb.Worksheets.Item("About")
$activateAboutResult = $wsAbout.Activate()
